$d = $word.ActiveDocument

# 1. Update the letter date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line into a street line and a city/state/zip
#    line, each on its own paragraph.
$d.Content.Find.Execute("2565 Greenrock Road, Milpitas CA 95035", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2565 Greenrock Road^pMilpitas, CA 95035", 2)

# The paragraph split above leaves the new "Milpitas, CA 95035" run with the
# document's default font instead of inheriting the Arial/11pt formatting
# used throughout the letter, so fix that up explicitly.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Milpitas, CA 95035") {
        $para.Range.Font.Name = "Arial"
        $para.Range.Font.NameAscii = "Arial"
        $para.Range.Font.NameOther = "Arial"
        $para.Range.Font.NameBi = "Arial"
        $para.Range.Font.Size = 11
        $para.Range.Font.SizeBi = 11
        break
    }
}

# 3. Remove the now-redundant blank "NoSpacing" paragraph that used to
#    separate "Board of Directors" from the signature block.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt.EndsWith("Board of Directors") -and $i -lt $count) {
        $nextPara = $d.Paragraphs.Item($i + 1)
        $nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)
        if ($nextText -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
